$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 20
$ws.Range("H20").Value = 1697.75
$ws.Range("I20").Value = 1697.75
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 1697.75
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -1467.75
$ws.Range("N20").ClearContents()
# Row 35
$ws.Range("H35").Value = 1697.75
$ws.Range("I35").Value = 1697.75
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 1697.75
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -1318.75
$ws.Range("N35").ClearContents()
# Row 51
$ws.Range("H51").Value = 2817.5518
$ws.Range("I51").Value = 2317.5
$ws.Range("J51").Value = 2948
$ws.Range("K51").Value = 2317.5
$ws.Range("L51").Value = 2948
$ws.Range("M51").Value = -1833.5
$ws.Range("N51").Value = -3916
# Row 74
$ws.Range("H74").Value = 2047.2693
$ws.Range("I74").Value = 1966.4783
$ws.Range("J74").Value = 2666.6667
$ws.Range("K74").Value = 1966.4783
$ws.Range("L74").Value = 2666.6667
$ws.Range("M74").Value = -1030.4783
$ws.Range("N74").Value = -4538.6667
# Row 77
$ws.Range("H77").Value = 2047.2693
$ws.Range("I77").Value = 1966.4783
$ws.Range("J77").Value = 2666.6667
$ws.Range("K77").Value = 9832.3915
$ws.Range("L77").Value = 13333.3335
$ws.Range("M77").Value = -5152.3915
$ws.Range("N77").Value = -22693.3335
# Row 92
$ws.Range("H92").Value = 541.8261
$ws.Range("I92").Value = 298.72223
$ws.Range("J92").Value = 1417
$ws.Range("K92").Value = 298.72223
$ws.Range("L92").Value = 1417
$ws.Range("M92").Value = 949.2777699999999
$ws.Range("N92").Value = -3913

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 97
$ws.Range("H97").Value = 715.4815
$ws.Range("I97").Value = 562.7619
$ws.Range("J97").Value = 1250
$ws.Range("K97").Value = 562.7619
$ws.Range("L97").Value = 1250
$ws.Range("M97").Value = -66.76189999999997
$ws.Range("N97").Value = -2242

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 61
$ws.Range("H61").Value = 9000
$ws.Range("J61").Value = 9000
$ws.Range("L61").Value = 9000
$ws.Range("N61").Value = -9626
# Row 94
$ws.Range("H94").Value = 4462.778
$ws.Range("I94").Value = 466.94116
$ws.Range("J94").Value = 11255.7
$ws.Range("K94").Value = 466.94116
$ws.Range("L94").Value = 11255.7
$ws.Range("M94").Value = -15.94116000000002
$ws.Range("N94").Value = -12157.7

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 47
$ws.Range("H47").Value = 6333.3335
$ws.Range("I47").Value = 1000
$ws.Range("K47").Value = 1000
$ws.Range("M47").Value = -434
# Row 94
$ws.Range("H94").Value = 1461.3334
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 1461.3334
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 1461.3334
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -2363.3334
# Row 99
$ws.Range("H99").Value = 1595
$ws.Range("I99").Value = 1430.4615
$ws.Range("K99").Value = 1430.4615
$ws.Range("M99").Value = 67.53850000000011
# Row 105
$ws.Range("H105").Value = 5650
$ws.Range("I105").Value = 1000
$ws.Range("K105").Value = 1000
$ws.Range("M105").Value = 747
# Row 122
$ws.Range("H122").Value = 1131.9259
$ws.Range("I122").Value = 1024.8572
$ws.Range("J122").Value = 1247.2307
$ws.Range("K122").Value = 3074.5716
$ws.Range("L122").Value = 3741.6921
$ws.Range("M122").Value = -624.5715999999998
$ws.Range("N122").Value = -8641.6921
# Row 126
$ws.Range("H126").Value = 1595
$ws.Range("I126").Value = 1430.4615
$ws.Range("K126").Value = 4291.3845
$ws.Range("M126").Value = -1821.3845

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 57
$ws.Range("H57").Value = 5700
$ws.Range("J57").Value = 11000
$ws.Range("L57").Value = 33000
$ws.Range("N57").Value = -34118
# Row 98
$ws.Range("H98").Value = 387.6875
$ws.Range("I98").Value = 300.23077
$ws.Range("J98").Value = 766.6667
$ws.Range("K98").Value = 900.69231
$ws.Range("L98").Value = 2300.0001
$ws.Range("M98").Value = 597.30769
$ws.Range("N98").Value = -5296.0001
# Row 131
$ws.Range("H131").Value = 822.1385
$ws.Range("I131").Value = 321.46155
$ws.Range("J131").Value = 947.3077
$ws.Range("K131").Value = 964.38465
$ws.Range("L131").Value = 2841.9231
$ws.Range("M131").Value = 4075.61535
$ws.Range("N131").Value = -12921.9231

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 4060.8
$ws.Range("I70").Value = 3944
$ws.Range("J70").Value = 4333.3335
$ws.Range("K70").Value = 3944
$ws.Range("L70").Value = 4333.3335
$ws.Range("M70").Value = -3674
$ws.Range("N70").Value = -4873.3335
# Row 73
$ws.Range("H73").Value = 4060.8
$ws.Range("I73").Value = 3944
$ws.Range("J73").Value = 4333.3335
$ws.Range("K73").Value = 3944
$ws.Range("L73").Value = 4333.3335
$ws.Range("M73").Value = -3008
$ws.Range("N73").Value = -6205.3335
# Row 122
$ws.Range("H122").Value = 4001.75
$ws.Range("I122").Value = 2002.3334
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 6007.0002
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -3557.0002
$ws.Range("N122").Value = -34900
# Row 126
$ws.Range("H126").Value = 1883.1428
$ws.Range("I126").Value = 1542
$ws.Range("J126").Value = 2338
$ws.Range("K126").Value = 4626
$ws.Range("L126").Value = 7014
$ws.Range("M126").Value = -2156
$ws.Range("N126").Value = -11954

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2060.1177
$ws.Range("I7").Value = 1909.3334
$ws.Range("J7").Value = 2422
$ws.Range("K7").Value = 1909.3334
$ws.Range("L7").Value = 2422
$ws.Range("M7").Value = -1797.3334
$ws.Range("N7").Value = -2646
# Row 122
$ws.Range("H122").Value = 62380.47
$ws.Range("I122").Value = 127801
$ws.Range("J122").Value = 4228.8887
$ws.Range("K122").Value = 383403
$ws.Range("L122").Value = 12686.6661
$ws.Range("M122").Value = -380953
$ws.Range("N122").Value = -17586.6661
# Row 126
$ws.Range("H126").Value = 2060.1177
$ws.Range("I126").Value = 1909.3334
$ws.Range("J126").Value = 2422
$ws.Range("K126").Value = 5728.0002
$ws.Range("L126").Value = 7266
$ws.Range("M126").Value = -3258.0002
$ws.Range("N126").Value = -12206

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 2632.2856
$ws.Range("I62").Value = 2265.7778
$ws.Range("J62").Value = 3292
$ws.Range("K62").Value = 2265.7778
$ws.Range("L62").Value = 3292
$ws.Range("M62").Value = -1641.7778
$ws.Range("N62").Value = -4540
# Row 65
$ws.Range("H65").Value = 2632.2856
$ws.Range("I65").Value = 2265.7778
$ws.Range("J65").Value = 3292
$ws.Range("K65").Value = 11328.889
$ws.Range("L65").Value = 16460
$ws.Range("M65").Value = -8208.888999999999
$ws.Range("N65").Value = -22700
# Row 136
$ws.Range("H136").Value = 3408.3845
$ws.Range("I136").Value = 3730.0293
$ws.Range("J136").Value = 2800.8333
$ws.Range("K136").Value = 11190.0879
$ws.Range("L136").Value = 8402.499899999999
$ws.Range("M136").Value = -8640.0879
$ws.Range("N136").Value = -13502.4999
